$d = $word.ActiveDocument
$d.Content.Find.Execute("Image Processing Toolbox, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
